$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "Virtual Machines" sheet: the Private DNS for VM01 (row 2) is
#    no longer known - change it to "N/A".
# ------------------------------------------------------------------
$vms = $wb.Worksheets.Item("Virtual Machines")
$vms.Range("E2").Value = "N/A"

# ------------------------------------------------------------------
# 2) Insert a new "AKS" worksheet right before "App Services" and
#    populate it with the AKS namespace/service/IP inventory.
# ------------------------------------------------------------------
$appServices = $wb.Worksheets.Item("App Services")
$aks = $wb.Worksheets.Add($appServices)
$aks.Name = "AKS"

# Header row
$aks.Cells.Item(1,1).Value = "Resource Group"
$aks.Cells.Item(1,2).Value = "AKS Server"
$aks.Cells.Item(1,3).Value = "Namespace"
$aks.Cells.Item(1,4).Value = "Service"
$aks.Cells.Item(1,5).Value = "Service IP"

# Data rows
$aks.Cells.Item(2,1).Value = "DefaultResourceGroup-WUS2"
$aks.Cells.Item(2,3).Value = "aks-istio-system"
$aks.Cells.Item(2,4).Value = "istiod-asm-1-19"
$aks.Cells.Item(2,5).Value = "10.0.98.233"

$aks.Cells.Item(3,1).Value = "DefaultResourceGroup-WUS2"
$aks.Cells.Item(3,3).Value = "calico-system"
$aks.Cells.Item(3,4).Value = "calico-kube-controllers-metrics"
$aks.Cells.Item(3,5).Value = "None"

$aks.Cells.Item(4,1).Value = "DefaultResourceGroup-WUS2"
$aks.Cells.Item(4,3).Value = "calico-system"
$aks.Cells.Item(4,4).Value = "calico-typha"
$aks.Cells.Item(4,5).Value = "10.0.203.206"

$aks.Cells.Item(5,1).Value = "DefaultResourceGroup-WUS2"
$aks.Cells.Item(5,3).Value = "default"
$aks.Cells.Item(5,4).Value = "kubernetes"
$aks.Cells.Item(5,5).Value = "10.0.0.1"

$aks.Cells.Item(6,1).Value = "DefaultResourceGroup-WUS2"
$aks.Cells.Item(6,3).Value = "gatekeeper-system"
$aks.Cells.Item(6,4).Value = "gatekeeper-webhook-service"
$aks.Cells.Item(6,5).Value = "10.0.220.234"

$aks.Cells.Item(7,1).Value = "DefaultResourceGroup-WUS2"
$aks.Cells.Item(7,3).Value = "kube-system"
$aks.Cells.Item(7,4).Value = "azure-policy-webhook-service"
$aks.Cells.Item(7,5).Value = "10.0.17.59"

$aks.Cells.Item(8,1).Value = "DefaultResourceGroup-WUS2"
$aks.Cells.Item(8,3).Value = "kube-system"
$aks.Cells.Item(8,4).Value = "kube-dns"
$aks.Cells.Item(8,5).Value = "10.0.0.10"

$aks.Cells.Item(9,1).Value = "DefaultResourceGroup-WUS2"
$aks.Cells.Item(9,3).Value = "kube-system"
$aks.Cells.Item(9,4).Value = "metrics-server"
$aks.Cells.Item(9,5).Value = "10.0.212.220"

# Column widths roughly matching the authored bestFit widths
$aks.Columns.Item(1).ColumnWidth = 25.81640625
$aks.Columns.Item(2).ColumnWidth = 9.6328125
$aks.Columns.Item(3).ColumnWidth = 16.6328125
$aks.Columns.Item(4).ColumnWidth = 26.81640625
$aks.Columns.Item(5).ColumnWidth = 11.453125

# ------------------------------------------------------------------
# 3) Reset the stale selection left on "Virtual Machines" (was G1)
#    and make the new "AKS" sheet the active / selected tab, with
#    column A selected (matches the authored selection state).
# ------------------------------------------------------------------
$vms.Range("A1").Select() | Out-Null
$sheet1 = $wb.Worksheets.Item("Sheet")
$sheet1.Range("A1").Select() | Out-Null
$appServices.Range("A1").Select() | Out-Null

$aks.Columns.Item(1).Select() | Out-Null
$aks.Activate() | Out-Null
